# Fix bug of user data duplication in worksheet
#
# The "Users" sheet accidentally contained duplicated / bogus rows
# (a second "Rita" row plus "sasha"/"dasha"/"Nadja" rows that were never
# meant to be there), and two of the legitimate rows had the wrong
# "Result" values. This cleans that up:
#   - Tanja's result (C3) is corrected from -1 to 9
#   - Rita's result (C6) is corrected from 1 to 6
#   - the four duplicated/bogus rows (7-10) are removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("C3").Value = 9
$ws.Range("C6").Value = 6

$ws.Rows("7:10").Delete()

$ws.Range("A7:XFD12").Select() | Out-Null

$wb.Save()
